$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# Cells whose target value is a numeric-looking string ("1","2","3") that must
# stay text (matching the source inlineStr cells) instead of being auto-coerced
# to a number by value-assignment type inference.
$textForceCells = @("A2", "A9", "A11")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$rowVals = @("3", "ANGEL MOLINA", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 21, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORRO", 94.71)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $rowVals[$i]
}

# Row 3
$rowVals = @("8-879-965", "CARLOS OREJUELA", "CABILLERO CALIFICADO", "Por horas", "No", "No", "No", 6.1, "16/01/2026", "30/01/2026", "31/01/2026", 87, 5, 38.12, 0, 8, 73.2, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 562.72)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $rowVals[$i]
}

# Row 4
$rowVals = @("8-970-1644", "DEBIN GONZALES", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 77, 3, 16.91, 0, 8, 54.12, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 368.69)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $rowVals[$i]
}

# Row 5
$rowVals = @("8-863-1584", "EDWIN FIGUEROA", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 79, 5, 28.19, 0, 8, 54.12, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 379.97)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $rowVals[$i]
}

# Row 6
$rowVals = @("8-921-1193", "EDWIN VEGA", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 77, 3, 16.91, 0, 8, 54.12, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 368.69)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "6").Value = $rowVals[$i]
}

# Row 7
$rowVals = @("1-716-753", "ESTEBAN PALACIO", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 84, 10, 56.37, 0, 8, 54.12, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 408.15)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "7").Value = $rowVals[$i]
}

# Row 8
$rowVals = @("8-952-1949", "JAIME MARTINEZ", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 29, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "0-441974958290", "BANCO GENERAL", "Ahorro", 130.79)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $rowVals[$i]
}

# Row 9
$rowVals = @("2", "JOSE BENITO CHIRINOS", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 45, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "", 202.95)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $rowVals[$i]
}

# Row 10
$rowVals = @("8-771-179", "JUAN RENTERIA", "ALBAÑIL CALIFICADO", "Por horas", "No", "No", "No", 6.1, "16/01/2026", "30/01/2026", "31/01/2026", 62, 4, 30.5, 0, 8, 73.2, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 408.7)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "10").Value = $rowVals[$i]
}

# Row 11
$rowVals = @("1", "LUIS FRIAS", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 15, 3, 16.91, 0, 8, 54.12, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORRO", 89.07)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "11").Value = $rowVals[$i]
}

# Row 12
$rowVals = @("8-730-847", "ORNELIS BARRIOS", "ALBAÑIL PRINPIPIANTE", "Por horas", "No", "No", "No", 5.09, "16/01/2026", "30/01/2026", "31/01/2026", 69, 3, 19.09, 0, 8, 61.08, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 375.39)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $rowVals[$i]
}

# Row 13
$rowVals = @("8-1042-173", "PEDRO CUEVAS", "AYUDANTE", "Por horas", "No", "No", "No", 4.51, "16/01/2026", "30/01/2026", "31/01/2026", 37, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, "0-472963882991", "BANCO GENERAL", "AHORRO", 166.87)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "13").Value = $rowVals[$i]
}

# Row 14
$rowVals = @("8-944-1234", "XAVIER SAENZ", "CABILLERO CALIFICADO", "Por horas", "No", "No", "No", 6.1, "16/01/2026", "30/01/2026", "31/01/2026", 79, 5, 38.12, 0, 8, 73.2, 0, 0, 0, 0, 0, 0, "", "BANCO GENERAL", "AHORROS", 513.92)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $rowVals[$i]
}
